$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 2791
$ws.Range("F5").Value = 7542
$ws.Range("F7").Value = 7743
$ws.Range("F9").Value = 41
$ws.Range("F10").Value = 24
$ws.Range("F11").Value = 6385
$ws.Range("F12").Value = 3310
$ws.Range("F15").Value = 33
$ws.Range("F16").Value = 30
$ws.Range("F17").Value = 31
$ws.Range("F19").Value = 19
$ws.Range("F21").Value = 6
$ws.Range("F22").Value = 298
$ws.Range("F24").Value = 3723
$ws.Range("F27").Value = 946
$ws.Range("F28").Value = 272
$ws.Range("F29").Value = 1381
$ws.Range("F31").Value = 37
$ws.Range("F32").Value = 2680
$ws.Range("F33").Value = 1666
$ws.Range("F34").Value = 25
$ws.Range("F36").Value = 41
$ws.Range("F37").Value = 3464
$ws.Range("F38").Value = 236
$ws.Range("F39").Value = 267
$ws.Range("F43").Value = 1342
$ws.Range("F44").Value = 237
$ws.Range("F46").Value = 613
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G4").Value = 108
$ws.Range("F8").Value = 36
$ws.Range("F12").Value = 34
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 127
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 127
$ws.Range("F5").Value = 2791
$ws.Range("G7").Value = 108
$ws.Range("F9").Value = 7542
$ws.Range("F10").Value = 7743
$ws.Range("F12").Value = 24
$ws.Range("F13").Value = 6385
$ws.Range("F14").Value = 3310
$ws.Range("F17").Value = 30
$ws.Range("F19").Value = 19
$ws.Range("F21").Value = 36
$ws.Range("F22").Value = 298
$ws.Range("F24").Value = 3723
$ws.Range("F27").Value = 34
$ws.Range("F29").Value = 946
$ws.Range("F30").Value = 272
$ws.Range("F31").Value = 1381
$ws.Range("F33").Value = 37
$ws.Range("F34").Value = 2680
$ws.Range("F35").Value = 1666
$ws.Range("F36").Value = 25
$ws.Range("F38").Value = 41
$ws.Range("F40").Value = 3464
$ws.Range("F41").Value = 236
$ws.Range("F42").Value = 267
$ws.Range("F46").Value = 1342
$ws.Range("F47").Value = 237
$ws.Range("F49").Value = 613
